$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 143
$ws.Cells.Item(143, 1).Value = 3683
$ws.Cells.Item(143, 2).Value = 87015
$ws.Cells.Item(143, 3).Value = 22889
$ws.Cells.Item(143, 4).Value = 6462
$ws.Cells.Item(143, 5).Value = 5013
$ws.Cells.Item(143, 6).Value = 2764
$ws.Cells.Item(143, 7).Value = 201801
$ws.Cells.Item(143, 8).Value = 2018
$ws.Cells.Item(143, 9).Value = 1
$ws.Cells.Item(143, 10).Value = 5063
$ws.Cells.Item(143, 11).Value = 1973
$ws.Cells.Item(143, 12).Value = 3003
$ws.Cells.Item(143, 13).Value = 3491
$ws.Cells.Item(143, 14).Value = 3923
$ws.Cells.Item(143, 15).Value = 97631
$ws.Cells.Item(143, 16).Value = 85245
$ws.Cells.Item(143, 17).Value = 74822
$ws.Cells.Item(143, 18).Value = 77071
$ws.Cells.Item(143, 19).Value = 81627
$ws.Cells.Item(143, 20).Value = 20210
$ws.Cells.Item(143, 21).Value = 21111
$ws.Cells.Item(143, 22).Value = 17163
$ws.Cells.Item(143, 23).Value = 18750
$ws.Cells.Item(143, 24).Value = 21556
$ws.Cells.Item(143, 25).Value = 7581
$ws.Cells.Item(143, 26).Value = 7499
$ws.Cells.Item(143, 27).Value = 5601
$ws.Cells.Item(143, 28).Value = 6152
$ws.Cells.Item(143, 29).Value = 8508
$ws.Cells.Item(143, 30).Value = 4842
$ws.Cells.Item(143, 31).Value = 3606
$ws.Cells.Item(143, 32).Value = 2913
$ws.Cells.Item(143, 33).Value = 1294
$ws.Cells.Item(143, 34).Value = 4420
$ws.Cells.Item(143, 35).Value = 7596
$ws.Cells.Item(143, 36).Value = 4013
$ws.Cells.Item(143, 37).Value = 5386
$ws.Cells.Item(143, 38).Value = 4054
$ws.Cells.Item(143, 39).Value = 4142

# Row 144
$ws.Cells.Item(144, 1).Value = 3274
$ws.Cells.Item(144, 2).Value = 80016
$ws.Cells.Item(144, 3).Value = 21051
$ws.Cells.Item(144, 4).Value = 6513
$ws.Cells.Item(144, 5).Value = 3637
$ws.Cells.Item(144, 6).Value = 2665
$ws.Cells.Item(144, 7).Value = 201802
$ws.Cells.Item(144, 8).Value = 2018
$ws.Cells.Item(144, 9).Value = 2
$ws.Cells.Item(144, 10).Value = 3683
$ws.Cells.Item(144, 11).Value = 5063
$ws.Cells.Item(144, 12).Value = 1973
$ws.Cells.Item(144, 13).Value = 4837
$ws.Cells.Item(144, 14).Value = 1640
$ws.Cells.Item(144, 15).Value = 87015
$ws.Cells.Item(144, 16).Value = 97631
$ws.Cells.Item(144, 17).Value = 85245
$ws.Cells.Item(144, 18).Value = 78286
$ws.Cells.Item(144, 19).Value = 49263
$ws.Cells.Item(144, 20).Value = 22889
$ws.Cells.Item(144, 21).Value = 20210
$ws.Cells.Item(144, 22).Value = 21111
$ws.Cells.Item(144, 23).Value = 15443
$ws.Cells.Item(144, 24).Value = 16832
$ws.Cells.Item(144, 25).Value = 6462
$ws.Cells.Item(144, 26).Value = 7581
$ws.Cells.Item(144, 27).Value = 7499
$ws.Cells.Item(144, 28).Value = 5895
$ws.Cells.Item(144, 29).Value = 4827
$ws.Cells.Item(144, 30).Value = 5013
$ws.Cells.Item(144, 31).Value = 4842
$ws.Cells.Item(144, 32).Value = 3606
$ws.Cells.Item(144, 33).Value = 4094
$ws.Cells.Item(144, 34).Value = 3179
$ws.Cells.Item(144, 35).Value = 2764
$ws.Cells.Item(144, 36).Value = 7596
$ws.Cells.Item(144, 37).Value = 4013
$ws.Cells.Item(144, 38).Value = 2712
$ws.Cells.Item(144, 39).Value = 1088

# Row 145
$ws.Cells.Item(145, 1).Value = 3158
$ws.Cells.Item(145, 2).Value = 79368
$ws.Cells.Item(145, 3).Value = 19995
$ws.Cells.Item(145, 4).Value = 7060
$ws.Cells.Item(145, 5).Value = 5276
$ws.Cells.Item(145, 6).Value = 4086
$ws.Cells.Item(145, 7).Value = 201803
$ws.Cells.Item(145, 8).Value = 2018
$ws.Cells.Item(145, 9).Value = 3
$ws.Cells.Item(145, 10).Value = 3274
$ws.Cells.Item(145, 11).Value = 3683
$ws.Cells.Item(145, 12).Value = 5063
$ws.Cells.Item(145, 13).Value = 3226
$ws.Cells.Item(145, 14).Value = 3269
$ws.Cells.Item(145, 15).Value = 80016
$ws.Cells.Item(145, 16).Value = 87015
$ws.Cells.Item(145, 17).Value = 97631
$ws.Cells.Item(145, 18).Value = 76128
$ws.Cells.Item(145, 19).Value = 79018
$ws.Cells.Item(145, 20).Value = 21051
$ws.Cells.Item(145, 21).Value = 22889
$ws.Cells.Item(145, 22).Value = 20210
$ws.Cells.Item(145, 23).Value = 17303
$ws.Cells.Item(145, 24).Value = 20244
$ws.Cells.Item(145, 25).Value = 6513
$ws.Cells.Item(145, 26).Value = 6462
$ws.Cells.Item(145, 27).Value = 7581
$ws.Cells.Item(145, 28).Value = 5355
$ws.Cells.Item(145, 29).Value = 8144
$ws.Cells.Item(145, 30).Value = 3637
$ws.Cells.Item(145, 31).Value = 5013
$ws.Cells.Item(145, 32).Value = 4842
$ws.Cells.Item(145, 33).Value = 3604
$ws.Cells.Item(145, 34).Value = 3976
$ws.Cells.Item(145, 35).Value = 2665
$ws.Cells.Item(145, 36).Value = 2764
$ws.Cells.Item(145, 37).Value = 7596
$ws.Cells.Item(145, 38).Value = 2127
$ws.Cells.Item(145, 39).Value = 2408

